$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New player roster data (Name, Position(s), Team), in the order Excel
# would have encountered them (first occurrence order drives the shared
# string table order).
$data = @(
    @("Shai Gilgeous-Alexander", "PG,SG",       "Oklahoma City Thunder"),
    @("Terry Rozier",            "PG",          "Miami Heat"),
    @("Jalen Williams",          "SG,SF,PF,C",  "Oklahoma City Thunder"),
    @("Jimmy Butler",            "SF,PF",       "Miami Heat"),
    @("Jordan Poole",            "PG,SG",       "Washington Wizards"),
    @("Dennis Schröder",         "PG",          "Golden State Warriors"),
    @("Kevin Porter Jr.",        "PG",          "LA Clippers"),
    @("Lauri Markkanen",         "SF,PF",       "Utah Jazz"),
    @("Zach LaVine",             "SG,SF",       "Chicago Bulls"),
    @("John Collins",            "PF,C",        "Utah Jazz"),
    @("CJ McCollum",             "PG,SG",       "New Orleans Pelicans"),
    @("Keyonte George",          "PG,SG",       "Utah Jazz"),
    @("Jordan Clarkson",         "SG,SF",       "Utah Jazz"),
    @("Tobias Harris",           "SF,PF",       "Detroit Pistons"),
    @("Kyrie Irving",            "PG,SG",       "Dallas Mavericks"),
    @("Joel Embiid",             "C",           "Philadelphia 76ers"),
    @("RJ Barrett",              "SG,SF,PF",    "Toronto Raptors")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
